$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menus_and_forms")

$ws.Range("E1").Value = "image_en"
$ws.Range("F1").Value = "audio_en"
$ws.Range("G1").Value = "image_fra"
$ws.Range("H1").Value = "audio_fra"

$ws.Range("H2").Select()
